$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.190233977619532
$ws.Range("C2").Value = 0.827873855544252
$ws.Range("D2").Value = 0.749745676500509
$ws.Range("E2").Value = 0.168870803662258
$ws.Range("F2").Value = 0.766022380467955
$ws.Range("G2").Value = 0.323499491353001
$ws.Range("H2").Value = 0.00712105798575788
$ws.Range("I2").Value = 0.0551373346897253
$ws.Range("J2").Value = 0.649643947100712
$ws.Range("K2").Value = 0.209562563580875
$ws.Range("L2").Value = 0.0354018311291963
$ws.Range("M2").Value = 0.838657171922686
$ws.Range("N2").Value = 0.0164801627670397
$ws.Range("O2").Value = 0.00488301119023398
$ws.Range("P2").Value = 0.00569684638860631
$ws.Range("Q2").Value = 0.491556459816887
$ws.Range("R2").Value = 0.750152594099695
$ws.Range("S2").Value = 0.0034587995930824
$ws.Range("T2").Value = 0.926144455747711
$ws.Range("U2").Value = 0.269379450661241
$ws.Range("V2").Value = 0.752390640895219
$ws.Range("W2").Value = 0.0474059003051882
$ws.Range("X2").Value = 0.133468972533062

$ws.Range("B3").Value = 0.0325534079348932
$ws.Range("C3").Value = 0.0939979654120041
$ws.Range("D3").Value = 0.0423194303153611
$ws.Range("E3").Value = 0.00244150559511699
$ws.Range("F3").Value = 0.196948118006104
$ws.Range("G3").Value = 0.0142421159715158
$ws.Range("H3").Value = 0.00162767039674466
$ws.Range("I3").Value = 0.623397761953205
$ws.Range("J3").Value = 0.134282807731434
$ws.Range("K3").Value = 0.0516785350966429
$ws.Range("L3").Value = 0.0551373346897253
$ws.Range("M3").Value = 0.0528992878942014
$ws.Range("N3").Value = 0.647812817904374
$ws.Range("O3").Value = 0.027263479145473
$ws.Range("P3").Value = 0.980061037639878
$ws.Range("Q3").Value = 0.15910478128179
$ws.Range("R3").Value = 0.0107833163784334
$ws.Range("S3").Value = 0.978026449643947
$ws.Range("T3").Value = 0.0034587995930824
$ws.Range("U3").Value = 0.0331637843336724
$ws.Range("V3").Value = 0.0166836215666328
$ws.Range("W3").Value = 0.0785350966429298
$ws.Range("X3").Value = 0.111495422177009

$ws.Range("B4").Value = 0.695625635808749
$ws.Range("C4").Value = 0.0703967446592065
$ws.Range("D4").Value = 0.201831129196338
$ws.Range("E4").Value = 0.17293997965412
$ws.Range("F4").Value = 0.00712105798575788
$ws.Range("G4").Value = 0.656154628687691
$ws.Range("H4").Value = 0.903763987792472
$ws.Range("I4").Value = 0.00447609359104781
$ws.Range("J4").Value = 0.209766022380468
$ws.Range("K4").Value = 0.0791454730417091
$ws.Range("L4").Value = 0.262665310274669
$ws.Range("M4").Value = 0.0170905391658189
$ws.Range("N4").Value = 0.071617497456765
$ws.Range("O4").Value = 0.0122075279755849
$ws.Range("P4").Value = 0.00264496439471007
$ws.Range("Q4").Value = 0.303967446592065
$ws.Range("R4").Value = 0.0138351983723296
$ws.Range("S4").Value = 0.017293997965412
$ws.Range("T4").Value = 0.0677517802644964
$ws.Range("U4").Value = 0.627263479145473
$ws.Range("V4").Value = 0.225432349949135
$ws.Range("W4").Value = 0.769074262461852
$ws.Range("X4").Value = 0.589827060020346

$ws.Range("B5").Value = 0.078942014242116
$ws.Range("C5").Value = 0.00366225839267548
$ws.Range("D5").Value = 0.00528992878942014
$ws.Range("E5").Value = 0.655137334689725
$ws.Range("F5").Value = 0.0295015259409969
$ws.Range("G5").Value = 0.00427263479145473
$ws.Range("H5").Value = 0.08646998982706
$ws.Range("I5").Value = 0.316581892166836
$ws.Range("J5").Value = 0.00488301119023398
$ws.Range("K5").Value = 0.658392675483215
$ws.Range("L5").Value = 0.644760935910478
$ws.Range("M5").Value = 0.0885045778229908
$ws.Range("N5").Value = 0.256154628687691
$ws.Range("O5").Value = 0.954018311291963
$ws.Range("P5").Value = 0.0115971515768057
$ws.Range("Q5").Value = 0.0400813835198372
$ws.Range("R5").Value = 0.222380467955239
$ws.Range("S5").Value = 0.00101729399796541
$ws.Range("T5").Value = 0.00101729399796541
$ws.Range("U5").Value = 0.0665310274669379
$ws.Range("V5").Value = 0.00427263479145473
$ws.Range("W5").Value = 0.100915564598169
$ws.Range("X5").Value = 0.163784333672431

